$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update item name in A3 (Mini Waterpomp x1 -> Mini Waterpomp x4 + buizen)
$ws.Range("A3").Value = "Mini Waterpomp x4 + buizen"

# Add new cost value in B3
$ws.Range("B3").Value = 21.64

# Update cost in C3 (4.37 -> 21.46)
$ws.Range("C3").Value = 21.46

# Update date in D3 to 11 April 2025
$ws.Range("D3").Value = (Get-Date -Year 2025 -Month 4 -Day 11 -Hour 0 -Minute 0 -Second 0)

# Column C should auto-fit its width like the diff shows (bestFit)
$ws.Columns.Item(3).AutoFit() | Out-Null

# Update the active selection to I10, matching the saved workbook view
$ws.Range("I10").Select() | Out-Null

$wb.Save()
